# Auto-generated edit script
# Updates "F" (want-to-go count) and one "G" (min price) cell per the
# commit "Update gh-pages to output generated at 456a3b4".
$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 344
$ws.Range("F3").Value = 3460
$ws.Range("F5").Value = 8099
$ws.Range("F7").Value = 64
$ws.Range("F8").Value = 2087
$ws.Range("F9").Value = 59
$ws.Range("F11").Value = 533
$ws.Range("F14").Value = 1052
$ws.Range("F18").Value = 1131
$ws.Range("F19").Value = 712
$ws.Range("F20").Value = 511
$ws.Range("F21").Value = 59
$ws.Range("F22").Value = 414
$ws.Range("F24").Value = 4560
$ws.Range("F26").Value = 47993
$ws.Range("F27").Value = 3954
$ws.Range("F28").Value = 24
$ws.Range("F30").Value = 718
$ws.Range("F31").Value = 16
$ws.Range("F32").Value = 63
$ws.Range("F33").Value = 823
$ws.Range("F35").Value = 558
$ws.Range("F38").Value = 560
$ws.Range("F40").Value = 954
$ws.Range("F41").Value = 111
$ws.Range("F42").Value = 149
$ws.Range("F43").Value = 1042
$ws.Range("F44").Value = 672
$ws.Range("F45").Value = 85
$ws.Range("F46").Value = 71
$ws.Range("F47").Value = 17
$ws.Range("F48").Value = 2444

# --- Sheet 2: 演出 ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("F4").Value = 235
$ws.Range("F11").Value = 108
$ws.Range("G11").Value = 180
$ws.Range("F14").Value = 33
$ws.Range("F15").Value = 77
$ws.Range("F19").Value = 7268
$ws.Range("F24").Value = 67
$ws.Range("F27").Value = 101
$ws.Range("F30").Value = 7

# --- Sheet 3: 本地生活 ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("F5").Value = 1470
$ws.Range("F7").Value = 622
$ws.Range("F8").Value = 2307
$ws.Range("F9").Value = 9227
$ws.Range("F10").Value = 1490
$ws.Range("F11").Value = 143
$ws.Range("F12").Value = 53

# --- Sheet 4: 全部类型 ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 3460
$ws.Range("F4").Value = 8099
$ws.Range("F5").Value = 1470
$ws.Range("F6").Value = 622
$ws.Range("F7").Value = 143
$ws.Range("F8").Value = 53
$ws.Range("F9").Value = 64
$ws.Range("F10").Value = 59
$ws.Range("F11").Value = 533
$ws.Range("F12").Value = 1052
$ws.Range("F13").Value = 235
$ws.Range("F15").Value = 1131
$ws.Range("F16").Value = 712
$ws.Range("F17").Value = 59
$ws.Range("F18").Value = 4560
$ws.Range("F20").Value = 108
$ws.Range("F23").Value = 3954
$ws.Range("F24").Value = 33
$ws.Range("F26").Value = 718
$ws.Range("F27").Value = 63
$ws.Range("F28").Value = 823
$ws.Range("F29").Value = 558
$ws.Range("F30").Value = 77
$ws.Range("F33").Value = 560
$ws.Range("F36").Value = 954
$ws.Range("F38").Value = 111
$ws.Range("F39").Value = 149
$ws.Range("F40").Value = 1042
$ws.Range("F41").Value = 672
$ws.Range("F43").Value = 85
$ws.Range("F45").Value = 71
$ws.Range("F46").Value = 17
$ws.Range("F48").Value = 2444

